$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the floating point precision of the existing A3 timestamp
$ws.Range("A3").Value = 44316.77361887153

# Append the new data row (row 4)
$ws.Range("A4").Value = 44317.77362215051
$ws.Range("B4").Value = 71645
$ws.Range("C4").Value = 60359
$ws.Range("D4").Value = 3044
$ws.Range("E4").Value = 1966
$ws.Range("F4").Value = 1396
$ws.Range("G4").Value = 18739
$ws.Range("H4").Value = 1283
$ws.Range("I4").Value = 792
$ws.Range("J4").Value = 189

# Apply the same date/time number format used for the rest of column A
# (style index 2 in the original workbook) to the new A4 cell.
$dateFormat = $ws.Range("A3").NumberFormat()
$ws.Range("A4").NumberFormat = $dateFormat
